$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Two withdrawals by "Puzziferri Domenico" on 25/05/2018 (serial 43245):
#  row 34: 2 Mt. of "Tela Pesante" (new material)
#  row 35: 2 Mt. of "Tela Leggera"

# Clone the per-column formatting used throughout the table (row 33 is the
# last existing data row) onto the two new rows before writing the values.
"A", "B", "C", "D", "E" | ForEach-Object {
    $col = $_
    $ws.Range("${col}33").Copy() | Out-Null
    $ws.Range("${col}34").PasteSpecial(-4122) | Out-Null
    $ws.Range("${col}35").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

$ws.Range("A34").Value = 43245
$ws.Range("B34").Value = "Puzziferri Domenico"
$ws.Range("C34").Value = "Tela Pesante"
$ws.Range("D34").Value = "Mt."
$ws.Range("E34").Value = 2

$ws.Range("A35").Value = 43245
$ws.Range("B35").Value = "Puzziferri Domenico"
$ws.Range("C35").Value = "Tela Leggera"
$ws.Range("D35").Value = "Mt."
$ws.Range("E35").Value = 2

$wb.Save()
